# The deck originally shipped with its slide-master theme (ppt/theme/theme1.xml,
# a "Red Violet" flavoured "Integral" theme) and its notes-master theme
# (ppt/theme/theme2.xml, the stock Office theme) swapped relative to what the
# author wanted: the slide master should use the stock "Office" colour
# palette. Re-point the twelve slide-master theme colour slots (dk1, lt1,
# dk2, lt2, accent1-6, hlink, folHlink) at the Office theme's RGB values so
# ppt/theme/theme1.xml ends up carrying the Office colour scheme, matching
# the target deck.

$p  = $ppt.ActivePresentation
$sm = $p.SlideMaster
$cs = $sm.ColorScheme

# Index -> (slot, target "Office" RGB hex)
#  1 dk1       000000
#  2 lt1       FFFFFF
#  3 dk2       44546A
#  4 lt2       E7E6E6
#  5 accent1   5B9BD5
#  6 accent2   ED7D31
#  7 accent3   A5A5A5
#  8 accent4   FFC000
#  9 accent5   4472C4
# 10 accent6   70AD47
# 11 hlink     0563C1
# 12 folHlink  954F72
$targetRgb = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

for ($i = 1; $i -le 12; $i++) {
    $color = $cs.Colors($i)
    $color.RGB = $targetRgb[$i - 1]
}
